# Commit: "added Json Writer class, change jsonUtil class to universal method"
#
# The underlying data edit (per the OOXML diff) swaps the "Математика" and
# "Физика" subject rows in the Statistics sheet: row 3 (previously Физика)
# becomes the Математика row, and row 5 (previously Математика) becomes the
# Физика row. Everything else (header row, Медицина row, Лингвистика row)
# stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current row 3 and row 5 values (read with Value2 - reliable getter)
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2

$a5 = $ws.Range("A5").Value2
$b5 = $ws.Range("B5").Value2
$c5 = $ws.Range("C5").Value2
$d5 = $ws.Range("D5").Value2
$e5 = $ws.Range("E5").Value2

# Write row 5's old content into row 3 ...
$ws.Range("A3").Value = $a5
$ws.Range("B3").Value = $b5
$ws.Range("C3").Value = $c5
$ws.Range("D3").Value = $d5
$ws.Range("E3").Value = $e5

# ... and row 3's old content into row 5 (full swap).
$ws.Range("A5").Value = $a3
$ws.Range("B5").Value = $b3
$ws.Range("C5").Value = $c3
$ws.Range("D5").Value = $d3
$ws.Range("E5").Value = $e3

# Re-setting the multi-line "university list" cells can make the host
# auto-apply an explicit custom row height; auto-fit both rows back so the
# sheet's row metadata is left as it was (no spurious ht/customHeight).
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(5).AutoFit()
